$d = $word.ActiveDocument

# --- Insert the two new blank paragraphs (done first, before bolding the
#     title, so the blank paragraph does not inherit bold formatting) ---

# Blank paragraph after "Use Case 7: Buch suchen"
$d.Paragraphs.Item(1).Range.InsertParagraphAfter()

# Blank paragraph after "Methode rufen:" (now paragraph 3)
$d.Paragraphs.Item(3).Range.InsertParagraphAfter()

# --- Make the title bold ---
$d.Paragraphs.Item(1).Range.Bold = 1

# --- Text fixes ---

# Method signature: "Buchsuche" -> "BuchSuche"
$d.Content.Find.Execute("gibListeZurBuchsucheZurueck", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "gibListeZurBuchSucheZurueck", 2) | Out-Null

# "Falls Titel oder autor oder isbn oder thema" -> "Falls Titel, autor, isbn oder thema"
$d.Content.Find.Execute("Falls Titel oder autor oder isbn oder thema übereinstimmen", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Falls Titel, autor, isbn oder thema übereinstimmen", 2) | Out-Null

# "Schliesse Datenbankverbindung" -> "Schließe Datenbankverbindung"
$d.Content.Find.Execute("Schliesse Datenbankverbindung", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Schließe Datenbankverbindung", 2) | Out-Null

# --- Restore the (Word-managed) _GoBack bookmark, now sitting right before
#     "isbn" in the "Falls Titel, autor, isbn oder thema ..." line ---
$fr = $d.Content
$fr.Find.ClearFormatting()
$fr.Find.Execute("isbn oder thema übereinstimmen") | Out-Null
$goBack = $d.Range($fr.Start, $fr.Start)
$d.Bookmarks.Add("_GoBack", $goBack) | Out-Null

Write-Output "done"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ("$i : [" + $p.Range.Text + "]")
}
